# Regione: REGIONE_CAMPANIA, Issuer: integrity:S1#VICAMPANIA
#
# The existing "Creazione_0" event row (row 2) is updated in place to
# become the newest "Aggiornamento_0" event, and the original
# "Creazione_0" event data is preserved by re-appending it as a new row
# (row 4), with an extra "Creazione_1" event inserted in between (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: turn the existing Creazione_0 row into the new Aggiornamento_0 event ---
# RDA (B) and PAZIENTE (C) stay the same for this patient/region, only the
# event id, workflow instance id, document id and event timestamp change.
$ws.Range("A2").Value = "Aggiornamento_0"
$ws.Range("D2").Value = "192989b8687539fd480210884fd676280d9149da24c19b00abe9298b11315bbb.d753ab6f72^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E2").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721298694310"
$ws.Range("F2").Value = "18-07-2024:12:31:41"

# --- Row 3: new Creazione_1 event ---
$ws.Range("A3").Value = "Creazione_1"
$ws.Range("B3").Value = "REGIONE_CAMPANIA"
$ws.Range("C3").Value = "NGNVCN92S19L259C^^^&2.16.840.1.113883.2.9.4.3.2&ISO"
$ws.Range("D3").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.80a82583af^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E3").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721298694310"
$ws.Range("F3").Value = "18-07-2024:12:31:35"

# --- Row 4: re-append the original Creazione_0 event (unchanged data) ---
$ws.Range("A4").Value = "Creazione_0"
$ws.Range("B4").Value = "REGIONE_CAMPANIA"
$ws.Range("C4").Value = "NGNVCN92S19L259C^^^&2.16.840.1.113883.2.9.4.3.2&ISO"
$ws.Range("D4").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.4d02802ff5^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E4").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721298575810"
$ws.Range("F4").Value = "18-07-2024:12:29:39"
